# planilhaDesafio.xlsx edit
# - Un-minimize the workbook window
# - Update the two generated usernames (vini001 -> vini007, lari001 -> lari007)
# - Move the active selection on the "register" sheet to C5
# - Resize columns A-C to fit the new (wider) username/password data
# - Set the print page setup to A4 portrait

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("register")

# Restore the application/workbook window to a normal (non-minimized) state.
$excel.WindowState = -4143
$wb.Windows.Item(1).WindowState = -4143

# Update the username values that changed.
$ws.Range("C2").Value = "vini007"
$ws.Range("C3").Value = "lari007"

# Widen the first three columns (no longer auto-fit) to better display the data.
$ws.Columns.Item(1).ColumnWidth = 11.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.333333333333334
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666

# Move/save the active cell selection.
[void]$ws.Range("C5").Select()

# Configure the page setup for printing (A4, portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
